# Apply the "想去人数" (want-to-go count) refresh + a couple of content
# corrections across the 展览 / 演出 / 本地生活 / 全部类型 sheets, matching
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------------
$wsExhibit.Range("F3").Value  = 2751
$wsExhibit.Range("F5").Value  = 944
$wsExhibit.Range("F7").Value  = 2406
$wsExhibit.Range("F8").Value  = 1860
$wsExhibit.Range("F9").Value  = 225
$wsExhibit.Range("F10").Value = 64
$wsExhibit.Range("F11").Value = 2516
$wsExhibit.Range("F14").Value = 59
$wsExhibit.Range("F16").Value = 133
$wsExhibit.Range("F17").Value = 123
$wsExhibit.Range("G17").Value = "已售罄"
$wsExhibit.Range("F18").Value = 9387
$wsExhibit.Range("F20").Value = 7305
$wsExhibit.Range("F21").Value = 11872
$wsExhibit.Range("F23").Value = 190
$wsExhibit.Range("F25").Value = 371
$wsExhibit.Range("F26").Value = 570
$wsExhibit.Range("F27").Value = 2658
$wsExhibit.Range("F30").Value = 2623
$wsExhibit.Range("F31").Value = 830
$wsExhibit.Range("F34").Value = 1003
$wsExhibit.Range("F35").Value = 357
$wsExhibit.Range("F36").Value = 50
$wsExhibit.Range("F37").Value = 551

# ---------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------
$wsShow.Range("F10").Value = 1187

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$wsLocal.Range("F4").Value = 169

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) -- combined listing
# ---------------------------------------------------------------------
$wsAll.Range("F5").Value  = 2751
$wsAll.Range("F8").Value  = 944
$wsAll.Range("F11").Value = 2406
$wsAll.Range("F13").Value = 1860
$wsAll.Range("F14").Value = 225
$wsAll.Range("F15").Value = 2516
$wsAll.Range("F19").Value = 59

# Row 20 becomes a brand-new listing ("海上钢琴师" concert); the former
# row 20 entry ("百鬼行代号鸢同人only") slides down into row 21, and the
# row that used to occupy row 21 ("第三届ICIC印象动漫节·宋媛媛专场签售劵")
# is dropped from this combined sheet.
$wsAll.Range("C20").Value = "杭州·海上钢琴师—一生必听的电影名曲《泰坦尼克号》《花样年华》《海上钢琴师》"
$wsAll.Range("D20").Value = "曙光路31号 浙江音乐厅"
$wsAll.Range("E20").Value = "2024.11.02 19:30-11.02 21:00"
$wsAll.Range("F20").Value = 6
$wsAll.Range("G20").Value = 100
$wsAll.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=90727"
$wsAll.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202408/J0vUXlhH1722908482422.png"

$wsAll.Range("C21").Value = "杭州·百鬼行代号鸢同人only"
$wsAll.Range("D21").Value = "保淑路2号 The Queen皇后"
$wsAll.Range("E21").Value = "2024.11.02 12:30-11.02 18:40"
$wsAll.Range("F21").Value = 133
$wsAll.Range("G21").Value = 168
$wsAll.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=92790"
$wsAll.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202409/bsDHN4VK1726910606937.jpeg"

$wsAll.Range("F22").Value = 9387
$wsAll.Range("F24").Value = 7305
$wsAll.Range("F25").Value = 11872
$wsAll.Range("F29").Value = 371
$wsAll.Range("F31").Value = 570
$wsAll.Range("F33").Value = 2658
$wsAll.Range("F46").Value = 551
